$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = 1.013391229189916
$ws.Range("C10").Value = 1.021155590884477
$ws.Range("C11").Value = 1.020276095821995
$ws.Range("C12").Value = 1.020146662069364
$ws.Range("C13").Value = 1.020421023121769
$ws.Range("C14").Value = 1.014783377783833
$ws.Range("C15").Value = 1.023574060176071
$ws.Range("C16").Value = 1.021529008095356
$ws.Range("C17").Value = 1.023084728102648
$ws.Range("C18").Value = 1.023181545086016
$ws.Range("C19").Value = 1.024900639787458
